$d = $word.ActiveDocument

$replacements = @(
    @("812×4=3248", "743×8=5944"),
    @("160×5=800", "242×2=484"),
    @("427×2=854", "591×2=1182"),
    @("528×8=4224", "344×3=1032"),
    @("847×6=5082", "607×6=3642"),
    @("915×9=8235", "183×6=1098"),
    @("859×2=1718", "629×2=1258"),
    @("586×9=5274", "610×5=3050"),
    @("804×5=4020", "889×7=6223"),
    @("533×2=1066", "976×9=8784"),
    @("339×6=2034", "625×7=4375"),
    @("137×5=685", "365×5=1825"),
    @("755×5=3775", "956×8=7648"),
    @("843×6=5058", "220×3=660"),
    @("144×8=1152", "119×9=1071"),
    @("887×5=4435", "819×8=6552"),
    @("127×2=254", "766×3=2298"),
    @("274×2=548", "466×7=3262"),
    @("988×5=4940", "602×2=1204"),
    @("577×9=5193", "125×3=375"),
    @("528×2=1056", "673×3=2019"),
    @("935×7=6545", "277×7=1939"),
    @("174×9=1566", "285×9=2565"),
    @("312×3=936", "151×7=1057"),
    @("195×5=975", "245×5=1225")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
